# Update gh-pages to output generated at 456a3b4
# Applies refreshed "want to go" / price counts across the four sheets of
# the workbook (展览, 演出, 本地生活, 全部类型).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 323
$ws1.Range("F3").Value  = 276
$ws1.Range("F5").Value  = 171
$ws1.Range("F6").Value  = 657
$ws1.Range("F8").Value  = 465
$ws1.Range("F9").Value  = 79
$ws1.Range("F10").Value = 505
$ws1.Range("F11").Value = 377
$ws1.Range("F12").Value = 63
$ws1.Range("F14").Value = 107
$ws1.Range("F15").Value = 194

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 42

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6217
$ws3.Range("F4").Value = 753
$ws3.Range("F5").Value = 1810
$ws3.Range("F6").Value = 122
$ws3.Range("G6").Value = "不可售"

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6217
$ws4.Range("F4").Value  = 753
$ws4.Range("F5").Value  = 1810
$ws4.Range("F6").Value  = 122
$ws4.Range("G6").Value  = "不可售"
$ws4.Range("F7").Value  = 323
$ws4.Range("F8").Value  = 276
$ws4.Range("F13").Value = 171
$ws4.Range("F16").Value = 657
$ws4.Range("F20").Value = 465
$ws4.Range("F22").Value = 79
$ws4.Range("F23").Value = 505
$ws4.Range("F25").Value = 377
$ws4.Range("F26").Value = 63
$ws4.Range("F30").Value = 107
$ws4.Range("F35").Value = 42
$ws4.Range("F36").Value = 194
